$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: set Index (A9) to 8
$ws.Cells.Item(9, 1).Value = 8

# Row 10: fill in new API entry
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "POST"
$ws.Cells.Item(10, 3).Value = "/api4/create_prescription/"
$ws.Cells.Item(10, 4).Value = "처방 데이터 추가"

# Update selection to C11 as seen in the diff
$ws.Range("C11").Select()
